$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "5Paanwala-Hara"
$ws.Range("B9").Value = "Njolssss"
